$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.297.31"
$ws.Range("E2").Value = "  -0.96%  "

$ws.Range("D3").Value = "1.561.68"
$ws.Range("E3").Value = "  -0.28%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.79%  "

$ws.Range("E6").Value = "  -0.69%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "44.31"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.51%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "23.59"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.59%  "

$ws.Range("E10").Value = "  -1.70%  "

$ws.Range("E11").Value = "  -1.03%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0894"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.96%  "

$ws.Range("D13").Value = "1.784.84"
$ws.Range("E13").Value = "  -0.22%  "

$ws.Range("D14").Value = "1.570.71"
$ws.Range("E14").Value = "  +0.32%  "

$ws.Range("D15").Value = "28.285.16"
$ws.Range("E15").Value = "  -0.96%  "

$ws.Range("E16").Value = "  -0.80%  "

$ws.Range("E17").Value = "  -2.07%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "60.95"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.91%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "227.59"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.47%  "

$ws.Range("E20").Value = "  -0.01%  "

$ws.Range("D21").Value = "0.0₃0675"
$ws.Range("E21").Value = "  -2.81%  "

$ws.Range("E22").Value = "  +0.04%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.92"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.31%  "

$ws.Range("E24").Value = "  -2.98%  "

$ws.Range("E25").Value = "  -2.44%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "150.20"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.53%  "

$ws.Range("E27").Value = "  -0.85%  "

$ws.Range("E28").Value = "  -0.54%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.32"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.93%  "

$ws.Range("E30").Value = "  +0.06%  "

$ws.Range("E31").Value = "  +2.23%  "

$ws.Range("E32").Value = "  -3.10%  "

$ws.Range("E33").Value = "  -1.29%  "

$ws.Range("E34").Value = "  -1.50%  "

$ws.Range("D35").Value = "1.376.31"
$ws.Range("E35").Value = "  -1.61%  "

$ws.Range("E36").Value = "  +1.45%  "

$ws.Range("E37").Value = "  -3.33%  "

$ws.Range("E38").Value = "  -0.30%  "

$ws.Range("E39").Value = "  +2.04%  "

$ws.Range("E40").Value = "  -2.09%  "

$ws.Range("E41").Value = "  -2.95%  "

$ws.Range("E42").Value = "  +2.84%  "

$ws.Range("E43").Value = "  +0.03%  "

$ws.Range("E44").Value = "  -0.50%  "

$ws.Range("E45").Value = "  -1.29%  "

$ws.Range("E46").Value = "  -3.62%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "62.03"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.14%  "

$ws.Range("E48").Value = "  -6.36%  "

$ws.Range("D49").Value = "1.697.55"
$ws.Range("E49").Value = "  -0.15%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "85.27"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.37%  "

$ws.Range("E51").Value = "  -2.24%  "
